# From v1.0.2 to v1.0.3
# The "Steps" content of TC3 and TC4 were swapped:
#   TC3 (rows 22-26) previously described the "atribuir/desatribuir" step,
#     now it describes the "realizar o empenho" step.
#   TC4 (rows 29-33) previously described the "realizar o empenho" step,
#     now it describes the "atribuir/desatribuir" step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tc3Action = $ws.Range("B26").Value()
$tc3Result = $ws.Range("D26").Value()
$tc4Action = $ws.Range("B33").Value()
$tc4Result = $ws.Range("D33").Value()

$ws.Range("B26").Value = $tc4Action
$ws.Range("D26").Value = $tc4Result

$ws.Range("B33").Value = $tc3Action
$ws.Range("D33").Value = $tc3Result
